# Auto-generated script applying Universalis market-price refresh values
# to the Diabolos_Profits workbook, per scheduled-runner commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 2811.8125
$ws.Range("I29").Value = 298.42856
$ws.Range("J29").Value = 4766.6665
$ws.Range("K29").Value = 895.28568
$ws.Range("L29").Value = 14299.9995
$ws.Range("M29").Value = -614.28568
$ws.Range("N29").Value = -14861.9995
$ws.Range("H58").Value = 3072.1428
$ws.Range("I58").Value = 294.42856
$ws.Range("J58").Value = 5849.857
$ws.Range("K58").Value = 883.28568
$ws.Range("L58").Value = 17549.571
$ws.Range("M58").Value = -733.28568
$ws.Range("N58").Value = -17849.571
$ws.Range("H74").Value = 3849.5454
$ws.Range("I74").Value = 3070
$ws.Range("J74").Value = 4499.1665
$ws.Range("K74").Value = 3070
$ws.Range("L74").Value = 4499.1665
$ws.Range("M74").Value = -2134
$ws.Range("N74").Value = -6371.1665
$ws.Range("H77").Value = 3849.5454
$ws.Range("I77").Value = 3070
$ws.Range("J77").Value = 4499.1665
$ws.Range("K77").Value = 15350
$ws.Range("L77").Value = 22495.8325
$ws.Range("M77").Value = -10670
$ws.Range("N77").Value = -31855.8325
$ws.Range("H138").Value = 7569.549
$ws.Range("I138").Value = 7146.846
$ws.Range("J138").Value = 7714.1577
$ws.Range("K138").Value = 21440.538
$ws.Range("L138").Value = 23142.4731
$ws.Range("M138").Value = -16300.538
$ws.Range("N138").Value = -33422.4731

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 7176.8335
$ws.Range("I5").Value = 765.25
$ws.Range("K5").Value = 765.25
$ws.Range("M5").Value = -653.25
$ws.Range("H32").Value = 160837.11
$ws.Range("J32").Value = 3156.4443
$ws.Range("L32").Value = 3156.4443
$ws.Range("N32").Value = -3730.4443
$ws.Range("H46").Value = 9638.799999999999
$ws.Range("I46").Value = 8000
$ws.Range("J46").Value = 10048.5
$ws.Range("K46").Value = 8000
$ws.Range("L46").Value = 10048.5
$ws.Range("M46").Value = -7681
$ws.Range("N46").Value = -10686.5
$ws.Range("H61").Value = 2028.619
$ws.Range("I61").Value = 1526.75
$ws.Range("K61").Value = 1526.75
$ws.Range("M61").Value = -1314.75
$ws.Range("H97").Value = 281.8846
$ws.Range("I97").Value = 298.66666
$ws.Range("J97").Value = 80.5
$ws.Range("K97").Value = 298.66666
$ws.Range("L97").Value = 80.5
$ws.Range("M97").Value = 197.33334
$ws.Range("N97").Value = -1072.5
$ws.Range("H110").Value = 35721172
$ws.Range("I110").Value = 38462800
$ws.Range("K110").Value = 38462800
$ws.Range("M110").Value = -38460755
$ws.Range("H132").Value = 50003930
$ws.Range("I132").Value = 55559364
$ws.Range("K132").Value = 166678092
$ws.Range("M132").Value = -166675562
$ws.Range("H136").Value = 2028.619
$ws.Range("I136").Value = 1526.75
$ws.Range("K136").Value = 4580.25
$ws.Range("M136").Value = -2030.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 7176.8335
$ws.Range("I4").Value = 765.25
$ws.Range("K4").Value = 765.25
$ws.Range("M4").Value = -650.25
$ws.Range("H94").Value = 9259514
$ws.Range("I94").Value = 10000115
$ws.Range("K94").Value = 10000115
$ws.Range("M94").Value = -9999664
$ws.Range("H99").Value = 669
$ws.Range("I99").Value = 654.6316
$ws.Range("J99").Value = 805.5
$ws.Range("K99").Value = 654.6316
$ws.Range("L99").Value = 805.5
$ws.Range("M99").Value = 843.3684
$ws.Range("N99").Value = -3801.5
$ws.Range("H134").Value = 2767.6453
$ws.Range("I134").Value = 2635.7144
$ws.Range("J134").Value = 3999
$ws.Range("K134").Value = 7907.1432
$ws.Range("L134").Value = 11997
$ws.Range("M134").Value = -5372.1432
$ws.Range("N134").Value = -17067

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 97
$ws.Range("I7").Value = 96.5
$ws.Range("K7").Value = 96.5
$ws.Range("M7").Value = 16.5
$ws.Range("H22").Value = 1707.1818
$ws.Range("I22").Value = 630.6667
$ws.Range("K22").Value = 630.6667
$ws.Range("M22").Value = -280.6667
$ws.Range("H31").Value = 2250.3333
$ws.Range("I31").Value = 861.7273
$ws.Range("K31").Value = 861.7273
$ws.Range("M31").Value = -566.7273
$ws.Range("H34").Value = 2250.3333
$ws.Range("I34").Value = 861.7273
$ws.Range("K34").Value = 861.7273
$ws.Range("M34").Value = -659.7273
$ws.Range("H86").Value = 5495.923
$ws.Range("J86").Value = 5789.8
$ws.Range("L86").Value = 5789.8
$ws.Range("N86").Value = -8035.8
$ws.Range("H89").Value = 5495.923
$ws.Range("J89").Value = 5789.8
$ws.Range("L89").Value = 28949
$ws.Range("N89").Value = -40181
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H132").Value = 836836.2
$ws.Range("I132").Value = 559025.3
$ws.Range("K132").Value = 1677075.9
$ws.Range("M132").Value = -1674545.9
$ws.Range("H134").Value = 3091.3333
$ws.Range("I134").Value = 1932.8334
$ws.Range("K134").Value = 5798.5002
$ws.Range("M134").Value = -3263.5002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 715.75
$ws.Range("I5").Value = 479.55554
$ws.Range("J5").Value = 909
$ws.Range("K5").Value = 1438.66662
$ws.Range("L5").Value = 2727
$ws.Range("M5").Value = -1326.66662
$ws.Range("N5").Value = -2951
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 100
$ws.Range("K23").Value = 300
$ws.Range("M23").Value = -65
$ws.Range("H56").Value = 6250
$ws.Range("I56").Value = 6250
$ws.Range("K56").Value = 6250
$ws.Range("M56").Value = -5720
$ws.Range("H68").Value = 3152.5334
$ws.Range("J68").Value = 3532.5
$ws.Range("L68").Value = 10597.5
$ws.Range("N68").Value = -12219.5
$ws.Range("H71").Value = 3152.5334
$ws.Range("J71").Value = 3532.5
$ws.Range("L71").Value = 31792.5
$ws.Range("N71").Value = -39904.5
$ws.Range("H107").Value = 2464.238
$ws.Range("J107").Value = 2527.7693
$ws.Range("L107").Value = 7583.3079
$ws.Range("N107").Value = -11423.3079
$ws.Range("H118").Value = 470.66666
$ws.Range("I118").Value = 468
$ws.Range("K118").Value = 1404
$ws.Range("M118").Value = -161
$ws.Range("H132").Value = 1684.8462
$ws.Range("J132").Value = 1900.5
$ws.Range("L132").Value = 17104.5
$ws.Range("N132").Value = -22164.5
$ws.Range("H135").Value = 715.75
$ws.Range("I135").Value = 479.55554
$ws.Range("J135").Value = 909
$ws.Range("K135").Value = 4315.99986
$ws.Range("L135").Value = 8181
$ws.Range("M135").Value = -1780.99986
$ws.Range("N135").Value = -13251

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H80").Value = 3369.6
$ws.Range("I80").Value = 1005
$ws.Range("K80").Value = 1005
$ws.Range("M80").Value = -7
$ws.Range("H83").Value = 3369.6
$ws.Range("I83").Value = 1005
$ws.Range("K83").Value = 5025
$ws.Range("M83").Value = -33
$ws.Range("H132").Value = 560887.8
$ws.Range("I132").Value = 1255724.8
$ws.Range("K132").Value = 3767174.4
$ws.Range("M132").Value = -3764644.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2875
$ws.Range("H22").Value = 6454413
$ws.Range("I22").Value = 1333.3334
$ws.Range("J22").Value = 16134032
$ws.Range("K22").Value = 1333.3334
$ws.Range("L22").Value = 16134032
$ws.Range("M22").Value = -1038.3334
$ws.Range("N22").Value = -16134622
$ws.Range("H27").Value = 6454413
$ws.Range("I27").Value = 1333.3334
$ws.Range("J27").Value = 16134032
$ws.Range("K27").Value = 1333.3334
$ws.Range("L27").Value = 16134032
$ws.Range("M27").Value = -1226.3334
$ws.Range("N27").Value = -16134246
$ws.Range("H46").Value = 3210.9033
$ws.Range("I46").Value = 3350
$ws.Range("J46").Value = 3177.52
$ws.Range("K46").Value = 3350
$ws.Range("L46").Value = 3177.52
$ws.Range("M46").Value = -3162
$ws.Range("N46").Value = -3553.52
$ws.Range("H55").Value = 922.5294
$ws.Range("I55").Value = 982.1429000000001
$ws.Range("K55").Value = 982.1429000000001
$ws.Range("M55").Value = -809.1429000000001
$ws.Range("H122").Value = 4416.6313
$ws.Range("I122").Value = 3665.1
$ws.Range("K122").Value = 10995.3
$ws.Range("M122").Value = -8545.299999999999
$ws.Range("H132").Value = 5619.207
$ws.Range("I132").Value = 3155.9473
$ws.Range("K132").Value = 9467.841899999999
$ws.Range("M132").Value = -6937.841899999999
$ws.Range("H136").Value = 6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2678.8235
$ws.Range("I122").Value = 2760.7693
$ws.Range("J122").Value = 2412.5
$ws.Range("K122").Value = 8282.3079
$ws.Range("L122").Value = 7237.5
$ws.Range("M122").Value = -5832.3079
$ws.Range("N122").Value = -12137.5
$ws.Range("H132").Value = 649902.0600000001
$ws.Range("I132").Value = 1254787.5
$ws.Range("J132").Value = 4690.933
$ws.Range("K132").Value = 3764362.5
$ws.Range("L132").Value = 14072.799
$ws.Range("M132").Value = -3761832.5
$ws.Range("N132").Value = -19132.799

Write-Host "Applied 238 cell updates across 8 sheets"
